$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in "generic" pair_kind (column J) for the practice rows 2-5 ---
$ws.Range("J2").Value = "generic"
$ws.Range("J3").Value = "generic"
$ws.Range("J4").Value = "generic"
$ws.Range("J5").Value = "generic"

# --- New "stim details" block appended near the bottom of the sheet ---
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

# Rows 29-36: month counts (A) paired with the media word_type (B)
$stimRows = @(
    @{ Row = 29; Month = 6; Type = "video" },
    @{ Row = 30; Month = 6; Type = "video" },
    @{ Row = 31; Month = 7; Type = "video" },
    @{ Row = 32; Month = 7; Type = "video" },
    @{ Row = 33; Month = 6; Type = "audio" },
    @{ Row = 34; Month = 6; Type = "audio" },
    @{ Row = 35; Month = 7; Type = "audio" },
    @{ Row = 36; Month = 7; Type = "audio" }
)

foreach ($item in $stimRows) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Month
    $ws.Cells.Item($item.Row, 2).Value = $item.Type
}
